$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label for column B was shortened (units dropped from the text).
$ws.Range("B1").Value = "Chlorophyll"

# Active selection moved back to B1 (top of the data column).
$ws.Range("B1").Select()
